# Limit No. Register Project Request feature
# Populate the request_list sheet with sample request rows.
# Rows are written in the same order the data was originally authored so the
# shared-string table is built up in the matching sequence (the "common"
# YCHERN/ASFLI/REGISTERPROJECT/PENDING request rows first, then the
# APPROVED rows, then the DEREGISTERPROJECT rows, then the CT113 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Request($row, $id, $fromUser, $toUser, $type, $status, $projectID) {
    $ws.Cells.Item($row, 2).Value = $fromUser
    $ws.Cells.Item($row, 3).Value = $toUser
    $ws.Cells.Item($row, 4).Value = $type
    $ws.Cells.Item($row, 5).Value = $status
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 6).Value = $projectID
}

Set-Request 6  5  "YCHERN" "ASFLI" "REGISTERPROJECT"   "PENDING"  1
Set-Request 7  6  "YCHERN" "ASFLI" "REGISTERPROJECT"   "PENDING"  2
Set-Request 8  7  "YCHERN" "ASFLI" "REGISTERPROJECT"   "PENDING"  3

Set-Request 2  1  "YCHERN" "ASFLI" "REGISTERPROJECT"   "APPROVED" 2
Set-Request 3  2  "YCHERN" "ASFLI" "REGISTERPROJECT"   "APPROVED" 1

Set-Request 4  3  "YCHERN" "ASFLI" "DEREGISTERPROJECT" "PENDING"  1
Set-Request 5  4  "YCHERN" "ASFLI" "DEREGISTERPROJECT" "PENDING"  2

Set-Request 9  8  "CT113"  "ASFLI" "REGISTERPROJECT"   "PENDING"  4
Set-Request 10 9  "CT113"  "ASFLI" "REGISTERPROJECT"   "PENDING"  5
Set-Request 11 10 "CT113"  "ASFLI" "REGISTERPROJECT"   "PENDING"  6
Set-Request 12 11 "CT113"  "ASFLI" "REGISTERPROJECT"   "PENDING"  7

# Match the source formatting: every populated data cell (including the
# blank trailing newTitle/newSupervisor columns) carries the default
# "Normal" style, same as the rest of the sheet.
$ws.Range("A2:H12").Style = "Normal"
